$wb = $excel.ActiveWorkbook

# Insert a new column before column N on the "Repayment Schedule" sheet,
# shifting the old N/O/P columns (Late / blank / Outstanding) one to the right.
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")
$wsSchedule.Columns("N:N").Insert()

# Make "Repayment Schedule" the active sheet (was "Transactions") and
# move its selection to R4. The "Transactions" sheet's own selection
# (E11) is unchanged by this edit.
$wsSchedule.Activate()
$wsSchedule.Range("R4").Select()
